$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Direct value updates (safe: Excel will not reinterpret these as numbers) ---
$ws.Range("D2").Value = '27.338.85'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '1.707.97'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("E6").Value = '  -1.65%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  -0.69%  '
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("E10").Value = '  -4.63%  '
$ws.Range("E11").Value = '  -0.57%  '
$ws.Range("E12").Value = '  -2.28%  '
$ws.Range("D13").Value = '1.942.45'
$ws.Range("E13").Value = '  -1.00%  '
$ws.Range("D14").Value = '1.705.71'
$ws.Range("E14").Value = '  -1.95%  '
$ws.Range("E15").Value = '  -0.79%  '
$ws.Range("D16").Value = '0.0₅8179'
$ws.Range("E16").Value = '  -1.60%  '
$ws.Range("E17").Value = '  -0.53%  '
$ws.Range("D18").Value = '27.330.45'
$ws.Range("E18").Value = '  -0.86%  '
$ws.Range("E19").Value = '  -2.80%  '
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("E21").Value = '  -2.24%  '
$ws.Range("E22").Value = '  -2.63%  '
$ws.Range("E23").Value = '  -1.67%  '
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("E25").Value = '  -2.97%  '
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("E27").Value = '  -2.65%  '
$ws.Range("E28").Value = '  -2.33%  '
$ws.Range("E29").Value = '  -2.39%  '
$ws.Range("E30").Value = '  -3.17%  '
$ws.Range("E31").Value = '  -0.79%  '
$ws.Range("E32").Value = '  -1.80%  '
$ws.Range("E33").Value = '  -1.00%  '
$ws.Range("E34").Value = '  -0.73%  '
$ws.Range("E35").Value = '  +1.40%  '
$ws.Range("E36").Value = '  -1.37%  '
$ws.Range("E37").Value = '  -2.07%  '
$ws.Range("E38").Value = '  -1.59%  '
$ws.Range("E39").Value = '  -0.49%  '
$ws.Range("E40").Value = '  -1.96%  '
$ws.Range("D41").Value = '1.046.69'
$ws.Range("E41").Value = '  -0.91%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("E42").Value = '  -1.05%  '
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("E44").Value = '  -0.55%  '
$ws.Range("D45").Value = '1.849.75'
$ws.Range("E45").Value = '  -0.94%  '
$ws.Range("E46").Value = '  +2.19%  '
$ws.Range("E47").Value = '  -2.10%  '
$ws.Range("E48").Value = '  +1.90%  '
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("E50").Value = '  -1.43%  '
$ws.Range("E51").Value = '  -0.86%  '

# --- Numeric-looking text updates: stage through a Text-formatted scratch cell
#     and PasteSpecial(values) so the destination keeps its original (default)
#     style/format and the value is preserved as text, not auto-converted to a Number. ---
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = '224.07'
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$scratch.Value = '0.5310'
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$scratch.Value = '0.2664'
$scratch.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$scratch.Value = '0.06618'
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$scratch.Value = '20.80'
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$scratch.Value = '0.07682'
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$scratch.Value = '4.514'
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$scratch.Value = '0.5824'
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$scratch.Value = '67.67'
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$scratch.Value = '215.01'
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$scratch.Value = '4.633'
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$scratch.Value = '10.42'
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$scratch.Value = '5.997'
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$scratch.Value = '143.91'
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$scratch.Value = '1.688'
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$scratch.Value = '0.1204'
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$scratch.Value = '7.241'
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$scratch.Value = '16.24'
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$scratch.Value = '0.05376'
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$scratch.Value = '1.293'
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$scratch.Value = '3.486'
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$scratch.Value = '3.425'
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$scratch.Value = '1.649'
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$scratch.Value = '2.863'
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$scratch.Value = '0.9508'
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$scratch.Value = '0.5860'
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$scratch.Value = '5.812'
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$scratch.Value = '0.8449'
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$scratch.Value = '1.003'
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$scratch.Value = '57.91'
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$scratch.Value = '0.4527'
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$scratch.Value = '1.006'
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$scratch.Value = '8.077'
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = $false
